# Typocrypha chatDatabase.xlsx - "Clarke's in" commit
# Inserts two new chat categories (bunbuku_1, wech_1) ahead of the existing
# doppelganger_1 category, and tweaks the doppelganger_1 dialogue line text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the doppelganger_1 dialogue line in place first (matches the
# author's edit order, and keeps the shared-string table layout identical).
$ws.Cells.Item(32, 2).Value = "Bear witness to my power!"

# Make room: push the existing doppelganger_1 block (rows 31-33) down by
# four rows so the two new categories can be inserted ahead of it.
$ws.Rows("31:34").Insert()

# Row 31-32: new "bunbuku_1" category
$ws.Cells.Item(31, 1).Value = "NEW_CATEGORY"
$ws.Cells.Item(31, 2).Value = "bunbuku_1"
$ws.Cells.Item(32, 1).Value = 100
$ws.Cells.Item(32, 2).Value = "I'm on fire, baby!"

# Row 33-34: new "wech_1" category
$ws.Cells.Item(33, 1).Value = "NEW_CATEGORY"
$ws.Cells.Item(33, 2).Value = "wech_1"
$ws.Cells.Item(34, 1).Value = 100
$ws.Cells.Item(34, 2).Value = "Please go away . . ."

# Restore the sheet selection recorded in the saved workbook.
$ws.Range("B25").Select()
